# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff:
#   - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#     (Overview!E2/F2, zh-cn!C2, de-de!C2)
#   - "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps bumped
#   - zh-cn / de-de summary columns narrowed to match the new (shorter) status text

$wb = $excel.ActiveWorkbook

$ovw  = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ovw.Range("E2").Value  = "Ready for handoff"
$ovw.Range("F2").Value  = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 (Latest HO Xliff Generate Date) and de-de!H2 (Latest Handoff Datetime)
# shared the same "2016-08-12 17:13:13" value; both move to "2016-08-12 17:14:08".
$ovw.Range("G2").Value  = "2016-08-12 17:14:08"
$dede.Range("H2").Value = "2016-08-12 17:14:08"

# zh-cn!H2 (Latest Handoff Datetime) moves from "2016-08-12 17:13:04" to "2016-08-12 17:13:56".
$zhcn.Range("H2").Value = "2016-08-12 17:13:56"

# --- Column widths ---
# The zh-cn/de-de columns re-sized to fit the new, shorter status text.
# (ColumnWidth is stored on a whole-pixel grid by the engine, so 16.3 is the
# input that lands on the closest attainable width to the target 17.2159881591797.)
$ovw.Columns.Item(5).ColumnWidth  = 16.3
$ovw.Columns.Item(6).ColumnWidth  = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
